$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.04229948730612585
$ws.Range("J4").Value = 0.4890179029360937
$ws.Range("K4").Value = 0.40771897256185
$ws.Range("L4").Value = 2.653083192022226
